$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 1. "Status" column text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (shared across Overview!E2:F2/E3:F3 and zh-cn!C2:C3 and de-de!C2:C3)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: fill in "Latest Target File" (I) and "Latest Handback File" (J)
#    for rows 2 and 3, and widen the relevant columns.
# ---------------------------------------------------------------------------
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b22839c2858f163d53a69936ce865049cb346990/e2e/1c1a1429-608f-4eb0-895f-769a9ab3dec4.md", "", "", "1c1a1429-608f-4eb0-895f-769a9ab3dec4.md")
$wsZh.Range("J2").Value = "1c1a1429-608f-4eb0-895f-769a9ab3dec4.cf44c2bac1fa52f7bde1877b6207e3bbfc8157c9.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-04 02:32:34"

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b22839c2858f163d53a69936ce865049cb346990/e2e/5700ef28-21c8-471a-8eb0-511e2b8fb4d1.md", "", "", "5700ef28-21c8-471a-8eb0-511e2b8fb4d1.md")
$wsZh.Range("J3").Value = "5700ef28-21c8-471a-8eb0-511e2b8fb4d1.ea4001eae3b1e6ebef36b7c5e8484bf5b73e7566.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-04 02:32:34"

$wsZh.Range("C1").ColumnWidth = 29.9777047293527
$wsZh.Range("I1").ColumnWidth = 40
$wsZh.Range("J1").ColumnWidth = 40

# ---------------------------------------------------------------------------
# 3. de-de sheet: fill in "Latest Target File" (I), "Latest Handback File" (J)
#    and "Latest Handback DateTime" (K) for rows 2 and 3, and widen columns.
# ---------------------------------------------------------------------------
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b22839c2858f163d53a69936ce865049cb346990/e2e/1c1a1429-608f-4eb0-895f-769a9ab3dec4.md", "", "", "1c1a1429-608f-4eb0-895f-769a9ab3dec4.md")
$wsDe.Range("J2").Value = "1c1a1429-608f-4eb0-895f-769a9ab3dec4.cf44c2bac1fa52f7bde1877b6207e3bbfc8157c9.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-04 02:32:41"

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b22839c2858f163d53a69936ce865049cb346990/e2e/5700ef28-21c8-471a-8eb0-511e2b8fb4d1.md", "", "", "5700ef28-21c8-471a-8eb0-511e2b8fb4d1.md")
$wsDe.Range("J3").Value = "5700ef28-21c8-471a-8eb0-511e2b8fb4d1.ea4001eae3b1e6ebef36b7c5e8484bf5b73e7566.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-04 02:32:41"

$wsDe.Range("C1").ColumnWidth = 29.9777047293527
$wsDe.Range("I1").ColumnWidth = 40
$wsDe.Range("J1").ColumnWidth = 40

Write-Host "Handback report generated"
